# "Latest changes added, final project"
#
# Adds a new "Mobile_Slider" worksheet as the last tab of the workbook,
# populates it with the mobile-slider table-header test data, updates the
# previously-active "Manage_Slider" sheet's remembered selection, and makes
# "Mobile_Slider" the new active/selected tab.

$wb = $excel.ActiveWorkbook

# Sheet that currently holds the selected/active tab (Manage_Slider).
$sliderSheet = $wb.Worksheets.Item("Manage_Slider")

# Create the new sheet by duplicating an existing simple 2-column sheet so it
# inherits the workbook's normal worksheet formatting/namespaces, then rename
# it and wipe its copied content.
$srcSheet = $wb.Worksheets.Item("Manage_News")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$srcSheet.Copy($null, $lastSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Mobile_Slider"
$newSheet.Cells.Clear()

# --- Populate the Mobile_Slider test data table ---
$newSheet.Range("A1").Value = "Table headers"

$newSheet.Range("A2").Value = "Image"
$newSheet.Range("B2").Value = "Status"
$newSheet.Range("C2").Value = "Action"

$newSheet.Range("A3").Value = "Edit button"
$newSheet.Range("B3").Value = "rgba(0, 123, 255, 1)"

$newSheet.Range("A4").Value = "Delete button"
$newSheet.Range("B4").Value = "rgba(220, 53, 69, 1)"

# --- Column widths for the new sheet ---
$newSheet.Columns.Item(1).ColumnWidth = 12
$newSheet.Columns.Item(2).ColumnWidth = 18.3

# --- Page setup (portrait) ---
$newSheet.PageSetup.Orientation = 1

# --- Restore Manage_Slider's (no longer active) remembered selection ---
$sliderSheet.Activate()
$sliderSheet.Range("B6").Select()

# --- Make Mobile_Slider the active/selected tab ---
$newSheet.Activate()
$newSheet.Range("L10").Select()
